$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.950.09"
$ws.Range("E2").Value = "  +2.63%  "

$ws.Range("D3").Value = "3.304.64"
$ws.Range("E3").Value = "  +2.34%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.20"
$ws.Range("E5").Value = "  +1.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.09"
$ws.Range("E6").Value = "  +2.60%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.602"
$ws.Range("E8").Value = "  +2.37%  "

$ws.Range("E9").Value = "  +4.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.68"
$ws.Range("E10").Value = "  -1.48%  "

$ws.Range("E11").Value = "  +2.79%  "

$ws.Range("D12").Value = "3.879.79"
$ws.Range("E12").Value = "  +2.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.137"
$ws.Range("E13").Value = "  +0.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.14"
$ws.Range("E14").Value = "  +4.73%  "

$ws.Range("D15").Value = "68.995.18"
$ws.Range("E15").Value = "  +2.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000173"
$ws.Range("E16").Value = "  +3.56%  "

$ws.Range("D17").Value = "3.323.28"
$ws.Range("E17").Value = "  +2.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.87"
$ws.Range("E18").Value = "  +1.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.67"
$ws.Range("E19").Value = "  +2.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "390.86"
$ws.Range("E20").Value = "  +4.72%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.78"
$ws.Range("E21").Value = "  +2.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.95"
$ws.Range("E22").Value = "  +1.48%  "

$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("E24").Value = "  +3.36%  "

$ws.Range("E25").Value = "  +2.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.188"
$ws.Range("E26").Value = "  +4.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.76"
$ws.Range("E27").Value = "  +2.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.80"
$ws.Range("E29").Value = "  +3.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.99"
$ws.Range("E30").Value = "  +1.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.16"
$ws.Range("E31").Value = "  +2.93%  "

$ws.Range("E32").Value = "  +4.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.19"
$ws.Range("E33").Value = "  +5.31%  "

$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("E35").Value = "  +4.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.28"
$ws.Range("E36").Value = "  +0.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.91"
$ws.Range("E37").Value = "  +3.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.838"
$ws.Range("E38").Value = "  -1.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.58"
$ws.Range("E39").Value = "  +0.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.60"
$ws.Range("E40").Value = "  +4.91%  "

$ws.Range("E41").Value = "  +1.12%  "

$ws.Range("E42").Value = "  -2.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.68"
$ws.Range("E43").Value = "  +2.83%  "

$ws.Range("E44").Value = "  +3.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.41"
$ws.Range("E45").Value = "  -0.23%  "

$ws.Range("D46").Value = "2.635.50"
$ws.Range("E46").Value = "  -2.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "341.86"
$ws.Range("E47").Value = "  -5.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0285"
$ws.Range("E48").Value = "  +2.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.26"
$ws.Range("E49").Value = "  +5.39%  "

$ws.Range("E50").Value = "  +0.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.31"
$ws.Range("E51").Value = "  +3.30%  "
